# Auto-generated: applies scheduled market-data refresh to Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2186413
$ws.Range("I137").Value = 4449350
$ws.Range("J137").Value = 1157805.5
$ws.Range("K137").Value = 13348050
$ws.Range("L137").Value = 3473416.5
$ws.Range("M137").Value = -13345500
$ws.Range("N137").Value = -3478516.5
$ws.Range("H138").Value = 2685.238
$ws.Range("I138").Value = 2045.9286
$ws.Range("J138").Value = 3963.8572
$ws.Range("K138").Value = 6137.7858
$ws.Range("L138").Value = 11891.5716
$ws.Range("M138").Value = -997.7857999999997
$ws.Range("N138").Value = -22171.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1116063.9
$ws.Range("I32").Value = 1194028.2
$ws.Range("J32").Value = 24562.6
$ws.Range("K32").Value = 1194028.2
$ws.Range("L32").Value = 24562.6
$ws.Range("M32").Value = -1193741.2
$ws.Range("N32").Value = -25136.6
$ws.Range("H74").Value = 8802058
$ws.Range("I74").Value = 5819702
$ws.Range("J74").Value = 18577558
$ws.Range("K74").Value = 5819702
$ws.Range("L74").Value = 18577558
$ws.Range("M74").Value = -5818828
$ws.Range("N74").Value = -18579306
$ws.Range("H77").Value = 8802058
$ws.Range("I77").Value = 5819702
$ws.Range("J77").Value = 18577558
$ws.Range("K77").Value = 29098510
$ws.Range("L77").Value = 92887790
$ws.Range("M77").Value = -29094142
$ws.Range("N77").Value = -92896526
$ws.Range("H132").Value = 34514.613
$ws.Range("I132").Value = 51587.05
$ws.Range("J132").Value = 3473.818
$ws.Range("K132").Value = 154761.15
$ws.Range("L132").Value = 10421.454
$ws.Range("M132").Value = -152231.15
$ws.Range("N132").Value = -15481.454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 9894739
$ws.Range("I99").Value = 3852084.5
$ws.Range("K99").Value = 3852084.5
$ws.Range("M99").Value = -3850586.5
$ws.Range("H105").Value = 46944.24
$ws.Range("J105").Value = 2085.9
$ws.Range("L105").Value = 2085.9
$ws.Range("N105").Value = -5579.9
$ws.Range("H134").Value = 2634.2666
$ws.Range("I134").Value = 2672.9666
$ws.Range("J134").Value = 2556.8667
$ws.Range("K134").Value = 8018.899800000001
$ws.Range("L134").Value = 7670.6001
$ws.Range("M134").Value = -5483.899800000001
$ws.Range("N134").Value = -12740.6001
$ws.Range("H137").Value = 52301.332
$ws.Range("J137").Value = 52301.332
$ws.Range("L137").Value = 52301.332
$ws.Range("N137").Value = -62501.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 959167.25
$ws.Range("I31").Value = 742.8
$ws.Range("J31").Value = 1821749.2
$ws.Range("K31").Value = 742.8
$ws.Range("L31").Value = 1821749.2
$ws.Range("M31").Value = -447.8
$ws.Range("N31").Value = -1822339.2
$ws.Range("H34").Value = 959167.25
$ws.Range("I34").Value = 742.8
$ws.Range("J34").Value = 1821749.2
$ws.Range("K34").Value = 742.8
$ws.Range("L34").Value = 1821749.2
$ws.Range("M34").Value = -540.8
$ws.Range("N34").Value = -1822153.2
$ws.Range("H140").Value = 28680
$ws.Range("J140").Value = 28680
$ws.Range("L140").Value = 28680
$ws.Range("N140").Value = -39040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1069.8813
$ws.Range("I68").Value = 562.3714
$ws.Range("J68").Value = 1810
$ws.Range("K68").Value = 1687.1142
$ws.Range("L68").Value = 5430
$ws.Range("M68").Value = -876.1142
$ws.Range("N68").Value = -7052
$ws.Range("H71").Value = 1069.8813
$ws.Range("I71").Value = 562.3714
$ws.Range("J71").Value = 1810
$ws.Range("K71").Value = 5061.3426
$ws.Range("L71").Value = 16290
$ws.Range("M71").Value = -1005.3426
$ws.Range("N71").Value = -24402
$ws.Range("H82").Value = 6783.1665
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 7939.8
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 23819.4
$ws.Range("M82").Value = -2594
$ws.Range("N82").Value = -24631.4
$ws.Range("H85").Value = 6783.1665
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 7939.8
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 23819.4
$ws.Range("M85").Value = -1596
$ws.Range("N85").Value = -26627.4
$ws.Range("H100").Value = 6354.8667
$ws.Range("J100").Value = 6521.2856
$ws.Range("L100").Value = 19563.8568
$ws.Range("N100").Value = -21185.8568
$ws.Range("H112").Value = 3816.6667
$ws.Range("I112").Value = 3666.6667
$ws.Range("J112").Value = 3966.6667
$ws.Range("K112").Value = 11000.0001
$ws.Range("L112").Value = 11900.0001
$ws.Range("M112").Value = -9892.000100000001
$ws.Range("N112").Value = -14116.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1418.5555
$ws.Range("J122").Value = 1370.875
$ws.Range("L122").Value = 4112.625
$ws.Range("N122").Value = -9012.625
$ws.Range("H126").Value = 2964.25
$ws.Range("I126").Value = 2735.6667
$ws.Range("J126").Value = 3650
$ws.Range("K126").Value = 8207.000100000001
$ws.Range("L126").Value = 10950
$ws.Range("M126").Value = -5737.000100000001
$ws.Range("N126").Value = -15890
$ws.Range("H132").Value = 2720933.2
$ws.Range("I132").Value = 4633049.5
$ws.Range("J132").Value = 3715.3684
$ws.Range("K132").Value = 13899148.5
$ws.Range("L132").Value = 11146.1052
$ws.Range("M132").Value = -13896618.5
$ws.Range("N132").Value = -16206.1052
$ws.Range("H135").Value = 32600
$ws.Range("J135").Value = 27000
$ws.Range("L135").Value = 27000
$ws.Range("N135").Value = -37140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1312.8235
$ws.Range("I93").Value = 1294.5333
$ws.Range("J93").Value = 1450
$ws.Range("K93").Value = 1294.5333
$ws.Range("L93").Value = 1450
$ws.Range("M93").Value = -46.53330000000005
$ws.Range("N93").Value = -3946
$ws.Range("H132").Value = 10110379
$ws.Range("I132").Value = 4588.357
$ws.Range("J132").Value = 17556752
$ws.Range("K132").Value = 13765.071
$ws.Range("L132").Value = 52670256
$ws.Range("M132").Value = -11235.071
$ws.Range("N132").Value = -52675316

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1738.9811
$ws.Range("I132").Value = 1169.325
$ws.Range("J132").Value = 3491.7693
$ws.Range("K132").Value = 3507.975
$ws.Range("L132").Value = 10475.3079
$ws.Range("M132").Value = -977.9750000000004
$ws.Range("N132").Value = -15535.3079

Write-Host "Applied 161 cell updates across 8 sheets"
